# Rotate the varying fields of rows 2-7 up by one (row2 -> row3, ..., row6 -> row7, row7 -> row2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry data which rotates between rows 2..7
$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R", "S")

# Capture the current (pre-edit) values for rows 2..7 for each relevant column
$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{}
    for ($r = 2; $r -le 7; $r++) {
        $orig[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# Apply the rotation: new row r (for r = 3..7) gets the old row (r-1) values;
# new row 2 gets the old row 7 values.
foreach ($col in $cols) {
    for ($r = 3; $r -le 7; $r++) {
        $ws.Range("$col$r").Value = $orig[$col][$r - 1]
    }
    $ws.Range("${col}2").Value = $orig[$col][7]
}
